$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.209.73'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '3.092.01'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.32'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.79%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '3.093.40'
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.456'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.43%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.35'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.106'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.398'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.30%  '
$ws.Range("D13").Value = '3.626.76'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("E14").Value = '  +0.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '25.35'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("E16").Value = '  -2.38%  '
$ws.Range("D17").Value = '57.307.37'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '3.096.43'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("E19").Value = '  -3.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.85'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.80%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '348.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("E26").Value = '  -1.78%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '0.0₃0870'
$ws.Range("E28").Value = '  -6.73%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.24'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("E31").Value = '  -0.58%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.88'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.89'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.43%  '
$ws.Range("E35").Value = '  -4.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.08'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.03'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.42%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.65'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.48%  '
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0660'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.53%  '
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("E42").Value = '  +0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.696'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").Value = '2.402.26'
$ws.Range("E44").Value = '  +5.75%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '36.72'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.31%  '
$ws.Range("E46").Value = '  +0.16%  '
$ws.Range("D47").Value = '3.130.96'
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0262'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  -4.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.90%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.77'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.98%  '
